# Update countries & provincias Spain - refresh COVID stats and fix sort order
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Last updated" timestamp (row 1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 9 de Agosto de 2020 a las 22:27"

# Refreshed per-country stats (and a few countries swap order after the refresh)
# Row 4
$ws.Cells.Item(4, 2).Value = 5181414
$ws.Cells.Item(4, 3).Value = 31691
$ws.Cells.Item(4, 4).Value = 2646388
$ws.Cells.Item(4, 5).Value = 2369595
$ws.Cells.Item(4, 7).Value = 361
$ws.Cells.Item(4, 8).Value = 165431
# Row 8
$ws.Cells.Item(8, 2).Value = 559859
$ws.Cells.Item(8, 3).Value = 6671
$ws.Cells.Item(8, 4).Value = 411474
$ws.Cells.Item(8, 5).Value = 137977
$ws.Cells.Item(8, 7).Value = 198
$ws.Cells.Item(8, 8).Value = 10408
# Row 22
$ws.Cells.Item(22, 2).Value = 217279
$ws.Cells.Item(22, 3).Value = 383
$ws.Cells.Item(22, 5).Value = 10618
# Row 31
$ws.Cells.Item(31, 4).Value = 78552
$ws.Cells.Item(31, 5).Value = 9985
# Row 34
$ws.Cells.Item(34, 2).Value = 83002
$ws.Cells.Item(34, 3).Value = 678
$ws.Cells.Item(34, 4).Value = 57533
$ws.Cells.Item(34, 5).Value = 24869
# Row 68
$ws.Cells.Item(68, 1).Value = "Costa Rica"
$ws.Cells.Item(68, 2).Value = 23286
$ws.Cells.Item(68, 3).Value = 484
$ws.Cells.Item(68, 4).Value = 7730
$ws.Cells.Item(68, 5).Value = 15321
$ws.Cells.Item(68, 7).Value = 7
$ws.Cells.Item(68, 8).Value = 235
# Row 69
$ws.Cells.Item(69, 1).Value = "Nepal"
$ws.Cells.Item(69, 2).Value = 22972
$ws.Cells.Item(69, 3).Value = 380
$ws.Cells.Item(69, 4).Value = 16353
$ws.Cells.Item(69, 5).Value = 6544
$ws.Cells.Item(69, 7).Value = 2
$ws.Cells.Item(69, 8).Value = 75
# Row 70
$ws.Cells.Item(70, 2).Value = 22818
$ws.Cells.Item(70, 3).Value = 565
$ws.Cells.Item(70, 4).Value = 10206
$ws.Cells.Item(70, 5).Value = 12205
$ws.Cells.Item(70, 7).Value = 17
$ws.Cells.Item(70, 8).Value = 407
# Row 76
$ws.Cells.Item(76, 2).Value = 16715
$ws.Cells.Item(76, 3).Value = 95
$ws.Cells.Item(76, 4).Value = 12926
$ws.Cells.Item(76, 5).Value = 3684
$ws.Cells.Item(76, 7).Value = 1
$ws.Cells.Item(76, 8).Value = 105
# Row 98
$ws.Cells.Item(98, 1).Value = "Mauritania"
$ws.Cells.Item(98, 2).Value = 6523
$ws.Cells.Item(98, 3).Value = 13
$ws.Cells.Item(98, 4).Value = 5527
$ws.Cells.Item(98, 5).Value = 839
$ws.Cells.Item(98, 8).Value = 157
# Row 99
$ws.Cells.Item(99, 1).Value = "Libano"
$ws.Cells.Item(99, 2).Value = 6517
$ws.Cells.Item(99, 3).Value = 294
$ws.Cells.Item(99, 4).Value = 2127
$ws.Cells.Item(99, 5).Value = 4312
$ws.Cells.Item(99, 8).Value = 78
# Row 103
$ws.Cells.Item(103, 2).Value = 5344
$ws.Cells.Item(103, 3).Value = 6
$ws.Cells.Item(103, 4).Value = 5106
$ws.Cells.Item(103, 5).Value = 179
# Row 109
$ws.Cells.Item(109, 1).Value = "Zimbabue"
$ws.Cells.Item(109, 2).Value = 4649
$ws.Cells.Item(109, 3).Value = 74
$ws.Cells.Item(109, 4).Value = 1437
$ws.Cells.Item(109, 5).Value = 3108
$ws.Cells.Item(109, 7).Value = 2
$ws.Cells.Item(109, 8).Value = 104
# Row 110
$ws.Cells.Item(110, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(110, 2).Value = 4641
$ws.Cells.Item(110, 4).Value = 1716
$ws.Cells.Item(110, 5).Value = 2866
$ws.Cells.Item(110, 8).Value = 59
# Row 120
$ws.Cells.Item(120, 1).Value = "Namibia"
$ws.Cells.Item(120, 2).Value = 2949
$ws.Cells.Item(120, 3).Value = 147
$ws.Cells.Item(120, 4).Value = 704
$ws.Cells.Item(120, 5).Value = 2226
$ws.Cells.Item(120, 7).Value = 3
$ws.Cells.Item(120, 8).Value = 19
# Row 121
$ws.Cells.Item(121, 1).Value = "Cabo Verde"
$ws.Cells.Item(121, 2).Value = 2858
$ws.Cells.Item(121, 3).Value = 23
$ws.Cells.Item(121, 4).Value = 2086
$ws.Cells.Item(121, 5).Value = 740
$ws.Cells.Item(121, 8).Value = 32
# Row 122
$ws.Cells.Item(122, 1).Value = "Sri Lanka"
$ws.Cells.Item(122, 2).Value = 2844
$ws.Cells.Item(122, 3).Value = 3
$ws.Cells.Item(122, 4).Value = 2579
$ws.Cells.Item(122, 5).Value = 254
$ws.Cells.Item(122, 8).Value = 11
# Row 131
$ws.Cells.Item(131, 2).Value = 2140
$ws.Cells.Item(131, 3).Value = 6
$ws.Cells.Item(131, 4).Value = 1346
$ws.Cells.Item(131, 5).Value = 787
$ws.Cells.Item(131, 7).Value = 1
$ws.Cells.Item(131, 8).Value = 7
# Row 136
$ws.Cells.Item(136, 2).Value = 1804
$ws.Cells.Item(136, 3).Value = 7
$ws.Cells.Item(136, 4).Value = 913
$ws.Cells.Item(136, 5).Value = 376
$ws.Cells.Item(136, 7).Value = 3
$ws.Cells.Item(136, 8).Value = 515
# Row 137
$ws.Cells.Item(137, 2).Value = 1697
$ws.Cells.Item(137, 3).Value = 19
$ws.Cells.Item(137, 4).Value = 1263
$ws.Cells.Item(137, 5).Value = 383
# Row 147
$ws.Cells.Item(147, 1).Value = "Siria"
$ws.Cells.Item(147, 2).Value = 1188
$ws.Cells.Item(147, 3).Value = 63
$ws.Cells.Item(147, 4).Value = 346
$ws.Cells.Item(147, 5).Value = 790
$ws.Cells.Item(147, 7).Value = 2
$ws.Cells.Item(147, 8).Value = 52
# Row 148
$ws.Cells.Item(148, 1).Value = "Burkina Faso"
$ws.Cells.Item(148, 2).Value = 1175
$ws.Cells.Item(148, 4).Value = 974
$ws.Cells.Item(148, 5).Value = 147
$ws.Cells.Item(148, 8).Value = 54
# Row 149
$ws.Cells.Item(149, 1).Value = "Niger"
$ws.Cells.Item(149, 2).Value = 1157
$ws.Cells.Item(149, 4).Value = 1057
$ws.Cells.Item(149, 5).Value = 31
$ws.Cells.Item(149, 8).Value = 69
# Row 152
$ws.Cells.Item(152, 2).Value = 1060
$ws.Cells.Item(152, 3).Value = 14
$ws.Cells.Item(152, 4).Value = 729
$ws.Cells.Item(152, 5).Value = 308
# Row 194
$ws.Cells.Item(194, 1).Value = "San Martin (Parte Francesa)"
$ws.Cells.Item(194, 2).Value = 78
$ws.Cells.Item(194, 3).Value = 25
$ws.Cells.Item(194, 4).Value = 44
$ws.Cells.Item(194, 5).Value = 31
$ws.Cells.Item(194, 8).Value = 3
# Row 195
$ws.Cells.Item(195, 1).Value = "Polinesia Francesa"
$ws.Cells.Item(195, 2).Value = 69
$ws.Cells.Item(195, 4).Value = 62
# Row 196
$ws.Cells.Item(196, 1).Value = "San Vicente y las Granadinas"
$ws.Cells.Item(196, 2).Value = 56
$ws.Cells.Item(196, 4).Value = 49
$ws.Cells.Item(196, 5).Value = 7
$ws.Cells.Item(196, 8).Value = 0
